$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 21 (2025Q3) metrics: total_customers, returning_customers, new_customers, recurrence_rate
$ws.Range("C21").Value = 278
$ws.Range("D21").Value = 241
$ws.Range("E21").Value = 37
$ws.Range("F21").Value = 69.05444126074498
